$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.580.83"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "1.817.04"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "228.12"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "0.559"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "34.72"
$ws.Range("E8").Value = "  +7.70%  "
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "2.077.52"
$ws.Range("E12").Value = "  +1.57%  "
$ws.Range("D13").Value = "11.36"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("D14").Value = "1.817.02"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "0.645"
$ws.Range("E15").Value = "  +3.27%  "
$ws.Range("D16").Value = "34.592.31"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "4.34"
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("D18").Value = "69.11"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").Value = "247.18"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "0.0₃0802"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  +5.55%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").Value = "172.66"
$ws.Range("E24").Value = "  +6.93%  "
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").Value = "7.45"
$ws.Range("E26").Value = "  +4.07%  "
$ws.Range("D27").Value = "16.76"
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("D28").Value = "0.117"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "4.05"
$ws.Range("E30").Value = "  +7.85%  "
$ws.Range("D31").Value = "0.0531"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").Value = "1.24"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.420.29"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "2.60"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").Value = "0.677"
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("D38").Value = "1.07"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").Value = "85.79"
$ws.Range("E40").Value = "  +4.92%  "
$ws.Range("D41").Value = "2.86"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("E42").Value = "  +3.74%  "
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").Value = "13.76"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("E46").Value = "  +3.09%  "
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").Value = "1.980.02"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").Value = "105.71"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("E51").Value = "  -0.02%  "
